$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows, per the "repull data,
# push all data, mean calculation" update.
$updates = @{
    8  = -1
    9  = -1
    10 = 0
    11 = -5
    12 = 4
    13 = 1
    14 = -1
    15 = -3
    16 = 1
    17 = 2
    18 = 7
    19 = -1
    20 = -4
    21 = -3
    22 = 2
    24 = 6
    25 = 3
    27 = -2
    28 = -4
    29 = 2
    31 = -1
    32 = 1
    33 = 3
    34 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
